# Clean up data columns to conform to specs: replace spaces with
# underscores in the "experimentDesign" (D) and "strain" (F) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2:D27 -> "Environmental Perturbation" becomes "Environmental_Perturbation"
$ws.Range("D2:D27").Value = "Environmental_Perturbation"

# F2:F27 -> "KN 99 alpha" becomes "KN 99_alpha"
$ws.Range("F2:F27").Value = "KN 99_alpha"

# Move the active selection from B2:B27 to F2:F27
$ws.Range("F2:F27").Select()
